$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the "Speaker Sex" side table (columns F:K) that is being removed
# as part of creating the OSF supplemental materials.
$ws.Range("F1:K21").ClearContents()

# Update the selection to match the new state of the sheet.
$ws.Range("E6:E7").Select()
